$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-03-21 Friday", $true, $true, $false, $false, $false, $true, 1, $false, "2025-03-22 Saturday", 2) | Out-Null
$d.Content.Find.Execute("35÷5=7, 0", $true, $true, $false, $false, $false, $true, 1, $false, "92÷3=30, 2", 2) | Out-Null
$d.Content.Find.Execute("99÷7=14, 1", $true, $true, $false, $false, $false, $true, 1, $false, "98÷7=14, 0", 2) | Out-Null
$d.Content.Find.Execute("35÷3=11, 2", $true, $true, $false, $false, $false, $true, 1, $false, "61÷2=30, 1", 2) | Out-Null
$d.Content.Find.Execute("16÷3=5, 1", $true, $true, $false, $false, $false, $true, 1, $false, "88÷4=22, 0", 2) | Out-Null
$d.Content.Find.Execute("98÷5=19, 3", $true, $true, $false, $false, $false, $true, 1, $false, "12÷4=3, 0", 2) | Out-Null
$d.Content.Find.Execute("15÷3=5, 0", $true, $true, $false, $false, $false, $true, 1, $false, "22÷3=7, 1", 2) | Out-Null
$d.Content.Find.Execute("65÷5=13, 0", $true, $true, $false, $false, $false, $true, 1, $false, "41÷6=6, 5", 2) | Out-Null
$d.Content.Find.Execute("68÷9=7, 5", $true, $true, $false, $false, $false, $true, 1, $false, "71÷4=17, 3", 2) | Out-Null
$d.Content.Find.Execute("94÷3=31, 1", $true, $true, $false, $false, $false, $true, 1, $false, "61÷4=15, 1", 2) | Out-Null
$d.Content.Find.Execute("97÷3=32, 1", $true, $true, $false, $false, $false, $true, 1, $false, "79÷9=8, 7", 2) | Out-Null
$d.Content.Find.Execute("12÷7=1, 5", $true, $true, $false, $false, $false, $true, 1, $false, "40÷4=10, 0", 2) | Out-Null
$d.Content.Find.Execute("19÷2=9, 1", $true, $true, $false, $false, $false, $true, 1, $false, "65÷7=9, 2", 2) | Out-Null
$d.Content.Find.Execute("48÷4=12, 0", $true, $true, $false, $false, $false, $true, 1, $false, "73÷7=10, 3", 2) | Out-Null
$d.Content.Find.Execute("46÷7=6, 4", $true, $true, $false, $false, $false, $true, 1, $false, "87÷4=21, 3", 2) | Out-Null
$d.Content.Find.Execute("79÷5=15, 4", $true, $true, $false, $false, $false, $true, 1, $false, "74÷5=14, 4", 2) | Out-Null
$d.Content.Find.Execute("42÷5=8, 2", $true, $true, $false, $false, $false, $true, 1, $false, "43÷2=21, 1", 2) | Out-Null
$d.Content.Find.Execute("45÷9=5, 0", $true, $true, $false, $false, $false, $true, 1, $false, "28÷9=3, 1", 2) | Out-Null
$d.Content.Find.Execute("95÷7=13, 4", $true, $true, $false, $false, $false, $true, 1, $false, "35÷7=5, 0", 2) | Out-Null
$d.Content.Find.Execute("29÷7=4, 1", $true, $true, $false, $false, $false, $true, 1, $false, "50÷9=5, 5", 2) | Out-Null
$d.Content.Find.Execute("28÷7=4, 0", $true, $true, $false, $false, $false, $true, 1, $false, "24÷2=12, 0", 2) | Out-Null
$d.Content.Find.Execute("11÷7=1, 4", $true, $true, $false, $false, $false, $true, 1, $false, "21÷9=2, 3", 2) | Out-Null
$d.Content.Find.Execute("58÷9=6, 4", $true, $true, $false, $false, $false, $true, 1, $false, "33÷9=3, 6", 2) | Out-Null
$d.Content.Find.Execute("31÷2=15, 1", $true, $true, $false, $false, $false, $true, 1, $false, "99÷6=16, 3", 2) | Out-Null
$d.Content.Find.Execute("53÷6=8, 5", $true, $true, $false, $false, $false, $true, 1, $false, "23÷8=2, 7", 2) | Out-Null
$d.Content.Find.Execute("40÷7=5, 5", $true, $true, $false, $false, $false, $true, 1, $false, "66÷9=7, 3", 2) | Out-Null
